# Automatische test-sync: 2025-08-18 20:42:50
# Append a new log row (row 6) to the "Logs" sheet, extend the conditional
# formatting ranges that depended on the last row, and refresh the
# "Dashboard" pivot-style summary count that changed as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New row of data (row 6)
$ws.Range("A6").Value = "Interne taak"
$ws.Range("B6").Value = "kwaliteit@testbedrijf123.nl"
$ws.Range("C6").Value = "Leg dit even neer bij Koen."
$ws.Range("D6").Value = "Onbekend"
$ws.Range("E6").Value = "Fout bij verwerken (forward_to_fallback() got an unexpected keyword argument 'afzender')"
$ws.Range("F6").Value = "2025-08-18 20:42:03"
$ws.Range("G6").Value = "Nee"
$ws.Range("H6").Value = "Ja"
$ws.Range("I6").Value = "Nee"
$ws.Range("J6").Value = "Nee"

# The conditional formatting sqref ranges for columns D, G, H, I, J covered
# rows 2:5 before; now that row 6 has data they must cover rows 2:6.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $ws.Range($col + "2:" + $col + "5")
    $newRange = $ws.Range($col + "2:" + $col + "6")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary: the new row's category is "Onbekend", whose
# tally goes from 1 to 2.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 2
